# Insert a new price record before the current row 229 ("Hortaliza, Macroferia
# Regional de Talca - Zapallo italiano" subconjunto), pushing all subsequent
# rows down by one (last row moves from 302 -> 303).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 229 (shifts rows 229:302 down to 230:303,
# carrying their formatting/styles with them, and grows the used range /
# dimension to A1:R303 automatically).
$ws.Rows.Item(229).Insert()

# Populate the newly inserted row 229 with the new record's data.
$ws.Range("A229").Value = 5
$ws.Range("B229").Value = "Macroferia Regional de Talca"
$ws.Range("C229").Value = "Maule"
$ws.Range("D229").Value = 44627
$ws.Range("E229").Value = 7
$ws.Range("F229").Value = 100112032
$ws.Range("G229").Value = "Zapallo italiano"
$ws.Range("H229").Value = "Sin especificar"
$ws.Range("I229").Value = "Primera"
$ws.Range("J229").Value = 400
$ws.Range("K229").Value = 7000
$ws.Range("L229").Value = 7000
$ws.Range("M229").Value = 7000
$ws.Range("N229").Value = "`$/caja 50 unidades"
$ws.Range("O229").Value = "Región del Maule"
$ws.Range("P229").Value = 140
$ws.Range("Q229").Value = 50
$ws.Range("R229").Value = "Hortaliza"
